# Scheduled data refresh: update cached market-price / profit figures
# (columns H..N = currentAveragePrice, currentAveragePriceNQ,
#  currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ,
#  LeveProfitHQ) across the per-job Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 10000
$ws.Cells.Item(18, 9).Value = 10000
$ws.Cells.Item(18, 11).Value = 10000
$ws.Cells.Item(18, 13).Value = -9716

$ws.Cells.Item(86, 8).Value = 11061.333
$ws.Cells.Item(86, 10).Value = 10815.833
$ws.Cells.Item(86, 12).Value = 10815.833
$ws.Cells.Item(86, 14).Value = -13061.833

$ws.Cells.Item(88, 8).Value = 3642.8572
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 3642.8572
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 3642.8572
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(88, 14).Value = -4454.8572

$ws.Cells.Item(89, 8).Value = 11061.333
$ws.Cells.Item(89, 10).Value = 10815.833
$ws.Cells.Item(89, 12).Value = 54079.165
$ws.Cells.Item(89, 14).Value = -65311.165

$ws.Cells.Item(91, 8).Value = 3642.8572
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 3642.8572
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 3642.8572
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(91, 14).Value = -6450.8572

$ws.Cells.Item(112, 8).Value = 4431.278
$ws.Cells.Item(112, 9).Value = 889.6667
$ws.Cells.Item(112, 10).Value = 4753.242
$ws.Cells.Item(112, 11).Value = 2669.0001
$ws.Cells.Item(112, 12).Value = 14259.726
$ws.Cells.Item(112, 13).Value = -1561.0001
$ws.Cells.Item(112, 14).Value = -16475.726

$ws.Cells.Item(125, 8).Value = 3254.162
$ws.Cells.Item(125, 9).Value = 1781
$ws.Cells.Item(125, 10).Value = 4052.125
$ws.Cells.Item(125, 11).Value = 16029
$ws.Cells.Item(125, 12).Value = 36469.125
$ws.Cells.Item(125, 13).Value = -13569
$ws.Cells.Item(125, 14).Value = -41389.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7352.4126
$ws.Cells.Item(32, 9).Value = 3842.3784
$ws.Cells.Item(32, 11).Value = 3842.3784
$ws.Cells.Item(32, 13).Value = -3555.3784

$ws.Cells.Item(33, 8).Value = 14857.143
$ws.Cells.Item(33, 9).Value = 15200
$ws.Cells.Item(33, 10).Value = 14000
$ws.Cells.Item(33, 11).Value = 15200
$ws.Cells.Item(33, 12).Value = 14000
$ws.Cells.Item(33, 13).Value = -14871
$ws.Cells.Item(33, 14).Value = -14658

$ws.Cells.Item(45, 8).Value = 38902.785
$ws.Cells.Item(45, 10).Value = 6842.4443
$ws.Cells.Item(45, 12).Value = 6842.4443
$ws.Cells.Item(45, 14).Value = -7596.4443

$ws.Cells.Item(61, 8).Value = 3139.3948
$ws.Cells.Item(61, 9).Value = 2975.4062
$ws.Cells.Item(61, 11).Value = 2975.4062
$ws.Cells.Item(61, 13).Value = -2763.4062

$ws.Cells.Item(74, 8).Value = 77046.375
$ws.Cells.Item(74, 9).Value = 83728.27
$ws.Cells.Item(74, 10).Value = 71392.46000000001
$ws.Cells.Item(74, 11).Value = 83728.27
$ws.Cells.Item(74, 12).Value = 71392.46000000001
$ws.Cells.Item(74, 13).Value = -82854.27
$ws.Cells.Item(74, 14).Value = -73140.46000000001

$ws.Cells.Item(77, 8).Value = 77046.375
$ws.Cells.Item(77, 9).Value = 83728.27
$ws.Cells.Item(77, 10).Value = 71392.46000000001
$ws.Cells.Item(77, 11).Value = 418641.35
$ws.Cells.Item(77, 12).Value = 356962.3
$ws.Cells.Item(77, 13).Value = -414273.35
$ws.Cells.Item(77, 14).Value = -365698.3

$ws.Cells.Item(118, 8).Value = 18410
$ws.Cells.Item(118, 10).Value = 18410
$ws.Cells.Item(118, 12).Value = 18410
$ws.Cells.Item(118, 14).Value = -21724

$ws.Cells.Item(132, 8).Value = 2981.303
$ws.Cells.Item(132, 9).Value = 2461.4814
$ws.Cells.Item(132, 10).Value = 5320.5
$ws.Cells.Item(132, 11).Value = 7384.4442
$ws.Cells.Item(132, 12).Value = 15961.5
$ws.Cells.Item(132, 13).Value = -4854.4442
$ws.Cells.Item(132, 14).Value = -21021.5

$ws.Cells.Item(136, 8).Value = 3139.3948
$ws.Cells.Item(136, 9).Value = 2975.4062
$ws.Cells.Item(136, 11).Value = 8926.2186
$ws.Cells.Item(136, 13).Value = -6376.2186

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(106, 8).Value = 110000
$ws.Cells.Item(106, 10).Value = 110000
$ws.Cells.Item(106, 12).Value = 110000
$ws.Cells.Item(106, 14).Value = -112524

$ws.Cells.Item(134, 8).Value = 3179.3333
$ws.Cells.Item(134, 9).Value = 1192.04
$ws.Cells.Item(134, 10).Value = 6728.0713
$ws.Cells.Item(134, 11).Value = 3576.12
$ws.Cells.Item(134, 12).Value = 20184.2139
$ws.Cells.Item(134, 13).Value = -1041.12
$ws.Cells.Item(134, 14).Value = -25254.2139

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(116, 8).Value = 4514.857
$ws.Cells.Item(116, 9).Value = 884.6667
$ws.Cells.Item(116, 10).Value = 7237.5
$ws.Cells.Item(116, 11).Value = 2654.0001
$ws.Cells.Item(116, 12).Value = 21712.5
$ws.Cells.Item(116, 13).Value = 787.9998999999998
$ws.Cells.Item(116, 14).Value = -28596.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 1789.125
$ws.Cells.Item(43, 9).Value = 1330.4286
$ws.Cells.Item(43, 11).Value = 1330.4286
$ws.Cells.Item(43, 13).Value = -1179.4286

$ws.Cells.Item(80, 8).Value = 62514936
$ws.Cells.Item(80, 9).Value = 125025950
$ws.Cells.Item(80, 10).Value = 3924
$ws.Cells.Item(80, 11).Value = 125025950
$ws.Cells.Item(80, 12).Value = 3924
$ws.Cells.Item(80, 13).Value = -125024952
$ws.Cells.Item(80, 14).Value = -5920

$ws.Cells.Item(83, 8).Value = 62514936
$ws.Cells.Item(83, 9).Value = 125025950
$ws.Cells.Item(83, 10).Value = 3924
$ws.Cells.Item(83, 11).Value = 625129750
$ws.Cells.Item(83, 12).Value = 19620
$ws.Cells.Item(83, 13).Value = -625124758
$ws.Cells.Item(83, 14).Value = -29604

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 13).ClearContents()

$ws.Cells.Item(61, 8).Value = 7997.4287
$ws.Cells.Item(61, 9).Value = 7339.0527
$ws.Cells.Item(61, 11).Value = 7339.0527
$ws.Cells.Item(61, 13).Value = -7137.0527

$ws.Cells.Item(113, 8).Value = 7997.4287
$ws.Cells.Item(113, 9).Value = 7339.0527
$ws.Cells.Item(113, 11).Value = 7339.0527
$ws.Cells.Item(113, 13).Value = -5169.0527

$ws.Cells.Item(136, 8).Value = 37465.27
$ws.Cells.Item(136, 9).Value = 57741.418
$ws.Cells.Item(136, 10).Value = 5728.696
$ws.Cells.Item(136, 11).Value = 173224.254
$ws.Cells.Item(136, 12).Value = 17186.088
$ws.Cells.Item(136, 13).Value = -170674.254
$ws.Cells.Item(136, 14).Value = -22286.088

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).ClearContents()

$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 921.4138
$ws.Cells.Item(113, 9).Value = 536.3333
$ws.Cells.Item(113, 11).Value = 1608.9999
$ws.Cells.Item(113, 13).Value = 561.0001

$ws.Cells.Item(122, 8).Value = 3613.682
$ws.Cells.Item(122, 9).Value = 2303.1428
$ws.Cells.Item(122, 10).Value = 5907.125
$ws.Cells.Item(122, 11).Value = 6909.428400000001
$ws.Cells.Item(122, 12).Value = 17721.375
$ws.Cells.Item(122, 13).Value = -4459.428400000001
$ws.Cells.Item(122, 14).Value = -22621.375
